# Branched from CRAN release
# Adds a new "2.5.0" benchmark row (models 3comp2 and sumclearances) to the
# httk-benchmarks data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Grow the table by one row - this also extends the table ref / autoFilter
# and the sheet dimension automatically.
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()
$newRowNum = $newRow.Range.Row

# Match formatting of the row above (left-aligned cell style) before filling
# in the new values.
$lastDataRow = $newRowNum - 1
$ws.Range("A" + $lastDataRow + ":R" + $lastDataRow).Copy()
$ws.Range("A" + $newRowNum + ":R" + $newRowNum).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A" + $newRowNum).Value = "2.5.0"
$ws.Range("B" + $newRowNum).Value = 1021
$ws.Range("C" + $newRowNum).Value = 1
$ws.Range("D" + $newRowNum).Value = 1
$ws.Range("E" + $newRowNum).Value = 0.9999
$ws.Range("F" + $newRowNum).Value = 0.9477
$ws.Range("G" + $newRowNum).Value = 353
$ws.Range("H" + $newRowNum).Value = 0.2716
$ws.Range("I" + $newRowNum).Value = 353
$ws.Range("J" + $newRowNum).Value = 1.508
$ws.Range("K" + $newRowNum).Value = 36
$ws.Range("L" + $newRowNum).Value = 0.9698
$ws.Range("M" + $newRowNum).Value = 80
$ws.Range("N" + $newRowNum).Value = 1.132
$ws.Range("O" + $newRowNum).Value = 80
$ws.Range("P" + $newRowNum).Value = 0.6466
$ws.Range("Q" + $newRowNum).Value = 863
$ws.Range("R" + $newRowNum).Value = "Added models 3comp2 and sumclearances"

# Move the view/selection onto the newly added row, like the author did.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 9
$ws.Range("R" + $newRowNum).Select()
